# Adds plots' underlying data for all five models (ResNet50, ResNet152, VGG16,
# VGG19, InceptionResNetV2, InceptionV3, EfficientNetB7): fills in missing
# Train/Validation timing columns (F/G/H) and Testing-time column (K) for the
# models that already had rows, inserts new rows for the newly-benchmarked
# models, and renames the "(tbc)" time-per-image header now that it has data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Time per image (tbc)" -> "Time per image (s)" ---
$ws.Range("H1").Value = "Time per image (s)"

# --- Row 2: ResNet50 (existing row, gains Train/Val time + Testing time data) ---
$ws.Range("A2").Value = "ResNet50"
$ws.Range("B2").Value = 0.9843
$ws.Range("C2").Value = 0.0504
$ws.Range("D2").Value = 0.95
$ws.Range("E2").Value = 0.2464
$ws.Range("F2").Value = 787
$ws.Range("G2").Value = 2700
$ws.Range("H2").Formula = "=F2/G2"
$ws.Range("I2").Value = 0.9333
$ws.Range("J2").Value = 0.037
$ws.Range("K2").Value = 39
$ws.Range("L2").Formula = "=36*5"
$ws.Range("M2").Formula = "=K2/L2"

# --- Row 3: ResNet152 (new row) ---
$ws.Range("A3").Value = "ResNet152"
$ws.Range("B3").Value = 0.9902
$ws.Range("C3").Value = 0.0214
$ws.Range("D3").Value = 0.9312
$ws.Range("E3").Value = 0.5332
$ws.Range("F3").Value = 2091
$ws.Range("G3").Value = 2700
$ws.Range("H3").Formula = "=F3/G3"
$ws.Range("I3").Value = 0.9556
$ws.Range("J3").Value = 0.251
$ws.Range("K3").Value = 92
$ws.Range("L3").Value = 180
$ws.Range("M3").Formula = "=K3/L3"

# --- Row 4: VGG16 (moved down from old row 3) ---
$ws.Range("A4").Value = "VGG16"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.013
$ws.Range("D4").Value = 0.975
$ws.Range("E4").Value = 0.084
$ws.Range("F4").Value = 2343
$ws.Range("G4").Value = 2700
$ws.Range("H4").Formula = "=F4/G4"
$ws.Range("I4").Value = 0.9778
$ws.Range("J4").Value = 0.0007
$ws.Range("K4").Value = 113
$ws.Range("L4").Value = 180
$ws.Range("M4").Formula = "=K4/L4"

# --- Row 5: VGG19 (new row) ---
$ws.Range("A5").Value = "VGG19"
$ws.Range("B5").Value = 0.9961
$ws.Range("C5").Value = 0.0189
$ws.Range("D5").Value = 0.9625
$ws.Range("E5").Value = 0.1609
$ws.Range("F5").Value = 2690
$ws.Range("G5").Value = 2700
$ws.Range("H5").Formula = "=F5/G5"
$ws.Range("I5").Value = 0.9833
$ws.Range("J5").Value = 0.0584
$ws.Range("K5").Value = 137
$ws.Range("L5").Value = 180
$ws.Range("M5").Formula = "=K5/L5"

# --- Row 6: InceptionResNetV2 (new row) ---
$ws.Range("A6").Value = "InceptionResNetV2"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0.0018
$ws.Range("D6").Value = 0.9375
$ws.Range("E6").Value = 0.8609
$ws.Range("F6").Value = 1196
$ws.Range("G6").Value = 2700
$ws.Range("H6").Formula = "=F6/G6"
$ws.Range("I6").Value = 0.95
$ws.Range("J6").Value = 0.5152
$ws.Range("K6").Value = 50
$ws.Range("L6").Value = 180
$ws.Range("M6").Formula = "=K6/L6"

# --- Row 7: InceptionV3 (moved down from old row 4) ---
$ws.Range("A7").Value = "InceptionV3"
$ws.Range("B7").Value = 0.9941
$ws.Range("C7").Value = 0.0878
$ws.Range("D7").Value = 0.9625
$ws.Range("E7").Value = 0.8078
$ws.Range("F7").Value = 474
$ws.Range("G7").Value = 2700
$ws.Range("H7").Formula = "=F7/G7"
$ws.Range("I7").Value = 0.9556
$ws.Range("J7").Value = 0.7297
$ws.Range("K7").Value = 24
$ws.Range("L7").Value = 180
$ws.Range("M7").Formula = "=K7/L7"

# --- Row 8: EfficientNetB7 (new row) ---
$ws.Range("A8").Value = "EfficientNetB7"
$ws.Range("B8").Value = 0.815
$ws.Range("C8").Value = 1.9407
$ws.Range("D8").Value = 0.8687
$ws.Range("E8").Value = 1.46
$ws.Range("F8").Value = 2363
$ws.Range("G8").Value = 2700
$ws.Range("H8").Formula = "=F8/G8"
$ws.Range("I8").Value = 0.9
$ws.Range("J8").Value = 0.8036
$ws.Range("K8").Value = 111
$ws.Range("L8").Value = 180
$ws.Range("M8").Formula = "=K8/L8"

# --- Column widths (cosmetic: column A widened for longer model names, new column M sized) ---
$ws.Columns.Item(1).ColumnWidth = 15.917
$ws.Columns.Item(13).ColumnWidth = 26.417

# --- View: scroll so column K is leftmost, select N7 (matches author's final on-screen state) ---
$ws.Range("N7").Select() | Out-Null
